# Automatische test-sync: 2025-08-06 20:04:50
# Appends a new logged mail-handling row to the "Logs" sheet and refreshes
# the "Dashboard" category summary (which is kept sorted by count, desc).

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Append the new row to the Logs sheet --------------------------------
$newRow = 13

$logs.Cells.Item($newRow, 1).Value = "Check jij even of dit nog geleverd kan worden?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #1: Check jij even of dit nog geleverd kan worden?"
$logs.Cells.Item($newRow, 4).Value = "Inkoop / Bestellingen"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-06 20:04:00"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# --- Extend the conditional formatting ranges to include the new row -----
# (Modify the existing rules' applies-to range in place so the rules
# themselves -- type/operator/dxfId/priority -- are left untouched.)
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$($col)2:$($col)12")
    $newRange = $logs.Range("$($col)2:$($col)13")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Refresh the Dashboard summary (sorted by Aantal, descending) --------
$dash.Cells.Item(2, 1).Value = "Inkoop / Bestellingen"
$dash.Cells.Item(2, 2).Value = 6
$dash.Cells.Item(3, 1).Value = "Planning / Afspraak"
$dash.Cells.Item(3, 2).Value = 5
